$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("DELETE, DROP e ALTER", $true, $false, $false, $false, $false, `
               $true, 1, $false, "DDL (DROP, ALTER) e Comando SHOW", 2)
